# Add periodic table element data (rows 55-84) to element_app sheet.
# Shared-string insertion order matches target workbook: all Symbols (col B),
# then all native-language names (col C), then English names (col A) last,
# so that Excel's shared string table de-dupes already-existing English names
# (Nitrogen, Oxygen, Fluorine, Phosphorus, Sulfur, Chlorine) correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: element symbols ---
$ws.Range("B55").Value = "H"
$ws.Range("B56").Value = "He"
$ws.Range("B57").Value = "Li"
$ws.Range("B58").Value = "Be"
$ws.Range("B59").Value = "B"
$ws.Range("B60").Value = "C"
$ws.Range("B61").Value = "N"
$ws.Range("B62").Value = "O"
$ws.Range("B63").Value = "F"
$ws.Range("B64").Value = "Ne"
$ws.Range("B65").Value = "Na"
$ws.Range("B66").Value = "Mg"
$ws.Range("B67").Value = "Al"
$ws.Range("B68").Value = "Si"
$ws.Range("B69").Value = "P"
$ws.Range("B70").Value = "S"
$ws.Range("B71").Value = "Cl"
$ws.Range("B72").Value = "Ar"
$ws.Range("B73").Value = "K"
$ws.Range("B74").Value = "Ca"
$ws.Range("B75").Value = "Sc"
$ws.Range("B76").Value = "Ti"
$ws.Range("B77").Value = "V"
$ws.Range("B78").Value = "Cr"
$ws.Range("B79").Value = "Mn"
$ws.Range("B80").Value = "Fe"
$ws.Range("B81").Value = "Co"
$ws.Range("B82").Value = "Ni"
$ws.Range("B83").Value = "Cu"
$ws.Range("B84").Value = "Zn"

# --- Column C: native-language element names ---
$ws.Range("C55").Value = "氫"
$ws.Range("C56").Value = "氦"
$ws.Range("C57").Value = "鋰"
$ws.Range("C58").Value = "鈹"
$ws.Range("C59").Value = "硼"
$ws.Range("C60").Value = "碳"
$ws.Range("C61").Value = "氮"
$ws.Range("C62").Value = "氧"
$ws.Range("C63").Value = "氟"
$ws.Range("C64").Value = "氖"
$ws.Range("C65").Value = "鈉"
$ws.Range("C66").Value = "鎂"
$ws.Range("C67").Value = "鋁"
$ws.Range("C68").Value = "矽（硅）"
$ws.Range("C69").Value = "磷"
$ws.Range("C70").Value = "硫"
$ws.Range("C71").Value = "氯"
$ws.Range("C72").Value = "氬"
$ws.Range("C73").Value = "鉀"
$ws.Range("C74").Value = "鈣"
$ws.Range("C75").Value = "鈧"
$ws.Range("C76").Value = "鈦"
$ws.Range("C77").Value = "釩"
$ws.Range("C78").Value = "鉻"
$ws.Range("C79").Value = "錳"
$ws.Range("C80").Value = "鐵"
$ws.Range("C81").Value = "鈷"
$ws.Range("C82").Value = "鎳"
$ws.Range("C83").Value = "銅"
$ws.Range("C84").Value = "鋅"

# --- Column A: English element names ---
$ws.Range("A55").Value = "Hydrogen"
$ws.Range("A56").Value = "Helium"
$ws.Range("A57").Value = "Lithium"
$ws.Range("A58").Value = "Beryllium"
$ws.Range("A59").Value = "Boron"
$ws.Range("A60").Value = "Carbon"
$ws.Range("A61").Value = "Nitrogen"
$ws.Range("A62").Value = "Oxygen"
$ws.Range("A63").Value = "Fluorine"
$ws.Range("A64").Value = "Neon"
$ws.Range("A65").Value = "Sodium"
$ws.Range("A66").Value = "Magnesium"
$ws.Range("A67").Value = "Aluminum"
$ws.Range("A68").Value = "Silicon"
$ws.Range("A69").Value = "Phosphorus"
$ws.Range("A70").Value = "Sulfur"
$ws.Range("A71").Value = "Chlorine"
$ws.Range("A72").Value = "Argon"
$ws.Range("A73").Value = "Potassium"
$ws.Range("A74").Value = "Calcium"
$ws.Range("A75").Value = "Scandium"
$ws.Range("A76").Value = "Titanium"
$ws.Range("A77").Value = "Vanadium"
$ws.Range("A78").Value = "Chromium"
$ws.Range("A79").Value = "Manganese"
$ws.Range("A80").Value = "Iron"
$ws.Range("A81").Value = "Cobalt"
$ws.Range("A82").Value = "Nickel"
$ws.Range("A83").Value = "Copper"
$ws.Range("A84").Value = "Zinc"

# --- Update the view: scroll down and select the newly added block ---
$ws.Range("A55:C84").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 48
$win.ScrollColumn = 1
